$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object "object[,]" 24,5
$inn = New-Object "object[,]" 24,6

$bf[0,0] = 1.02
$bf[0,1] = 1.014668797375967
$bf[0,2] = 1.021211119736525
$bf[0,3] = 0.9926147277508489
$bf[0,4] = 1.026858221286603
$inn[0,0] = 1.027611606393266
$inn[0,1] = 1.01989871326918
$inn[0,2] = 1.024049120349957
$inn[0,3] = 0.9955398523336033
$inn[0,4] = 1.029679641441378
$inn[0,5] = 1.010775762258241
$bf[1,0] = 1.02
$bf[1,1] = 1.015479718667332
$bf[1,2] = 1.021789817107072
$bf[1,3] = 0.9936372048519304
$bf[1,4] = 1.027926763525179
$inn[1,0] = 1.027746054806674
$inn[1,1] = 1.020345044725622
$inn[1,2] = 1.024434846105535
$inn[1,3] = 0.9963617723202692
$inn[1,4] = 1.030555080408866
$inn[1,5] = 1.010922686223704
$bf[2,0] = 1.02
$bf[2,1] = 1.016004570188592
$bf[2,2] = 1.022163994597182
$bf[2,3] = 0.9942998659930995
$bf[2,4] = 1.028618306622823
$inn[2,0] = 1.027831202631635
$inn[2,1] = 1.020633348066851
$inn[2,2] = 1.024683460171452
$inn[2,3] = 0.9968940712668345
$inn[2,4] = 1.031121083325464
$inn[2,5] = 1.011017582247846
$bf[3,0] = 1.02
$bf[3,1] = 1.016225247861233
$bf[3,2] = 1.022321230528245
$bf[3,3] = 0.9945786998346017
$bf[3,4] = 1.028909060491324
$inn[3,0] = 1.027866555275906
$inn[3,1] = 1.02075442936866
$inn[3,2] = 1.024787742636343
$inn[3,3] = 0.997117960005301
$inn[3,4] = 1.031358918909325
$inn[3,5] = 1.011057434626693
$bf[4,0] = 1.02
$bf[4,1] = 1.016262302339686
$bf[4,2] = 1.022347627080899
$bf[4,3] = 0.9946255319796338
$bf[4,4] = 1.028957881033288
$inn[4,0] = 1.027872465110583
$inn[4,1] = 1.020774752302138
$inn[4,2] = 1.024805238313151
$inn[4,3] = 0.9971555583673453
$inn[4,4] = 1.031398845956599
$inn[4,5] = 1.0110641235487
$bf[5,0] = 1.02
$bf[5,1] = 1.016007518777872
$bf[5,2] = 1.022166095861778
$bf[5,3] = 0.9943035907982488
$bf[5,4] = 1.028622191576635
$inn[5,0] = 1.027831676759041
$inn[5,1] = 1.020634966439524
$inn[5,2] = 1.024684854523276
$inn[5,3] = 0.9968970624462087
$inn[5,4] = 1.031124261737455
$inn[5,5] = 1.011018114922399
$bf[6,0] = 1.02
$bf[6,1] = 1.014942824059693
$bf[6,2] = 1.021406749887903
$bf[6,3] = 0.9929600610674301
$bf[6,4] = 1.027219313856296
$inn[6,0] = 1.027657426305596
$inn[6,1] = 1.02004965660706
$inn[6,2] = 1.024179679625867
$inn[6,3] = 0.995817528259106
$inn[6,4] = 1.029975595860747
$inn[6,5] = 1.010825451637862
$bf[7,0] = 1.02
$bf[7,1] = 1.013067758989844
$bf[7,2] = 1.020066624788036
$bf[7,3] = 0.9906006454969559
$bf[7,4] = 1.024748251722803
$inn[7,0] = 1.027336245436603
$inn[7,1] = 1.019014458327734
$inn[7,2] = 1.023282070643777
$inn[7,3] = 0.9939188001724441
$inn[7,4] = 1.027947976385528
$inn[7,5] = 1.010484639966816
$bf[8,0] = 1.02
$bf[8,1] = 1.011818507387199
$bf[8,2] = 1.019171918702394
$bf[8,3] = 0.989033133672735
$bf[8,4] = 1.023101580970847
$inn[8,0] = 1.027112670302164
$inn[8,1] = 1.018321825269327
$inn[8,2] = 1.022678739519785
$inn[8,3] = 0.9926553831429383
$inn[8,4] = 1.026593906057375
$inn[8,5] = 1.010256568805982
$bf[9,0] = 1.02
$bf[9,1] = 1.011277771307175
$bf[9,2] = 1.018784215288543
$bf[9,3] = 0.988355674866747
$bf[9,4] = 1.022388728876544
$inn[9,0] = 1.027013626605431
$inn[9,1] = 1.018021326540884
$inn[9,2] = 1.022416336359236
$inn[9,3] = 0.9921088820399291
$inn[9,4] = 1.026007036147011
$inn[9,5] = 1.010157610892357
$bf[10,0] = 1.02
$bf[10,1] = 1.011076948847768
$bf[10,2] = 1.01864016303293
$bf[10,3] = 0.9881042295826724
$bf[10,4] = 1.022123969494695
$inn[10,0] = 1.026976502411515
$inn[10,1] = 1.017909621138049
$inn[10,2] = 1.022318695442314
$inn[10,3] = 0.9919059725120875
$inn[10,4] = 1.025788964731411
$inn[10,5] = 1.010120823557347
$bf[11,0] = 1.02
$bf[11,1] = 1.011120024524571
$bf[11,2] = 1.018671064603439
$bf[11,3] = 0.9881581567098651
$bf[11,4] = 1.022180760113161
$inn[11,0] = 1.026984480823262
$inn[11,1] = 1.017933586238617
$inn[11,2] = 1.022339647547929
$inn[11,3] = 0.9919494934313052
$inn[11,4] = 1.025835745482229
$inn[11,5] = 1.010128715916085
$bf[12,0] = 1.02
$bf[12,1] = 1.011261170620582
$bf[12,2] = 1.018772308730667
$bf[12,3] = 0.9883348863814464
$bf[12,4] = 1.022366843257785
$inn[12,0] = 1.027010564736894
$inn[12,1] = 1.018012094702561
$inn[12,2] = 1.022408268847287
$inn[12,3] = 0.9920921077337197
$inn[12,4] = 1.025989011960969
$inn[12,5] = 1.010154570649051
$bf[13,0] = 1.02
$bf[13,1] = 1.011348139532986
$bf[13,2] = 1.018834683068171
$bf[13,3] = 0.9884438009545853
$bf[13,4] = 1.022481498635404
$inn[13,0] = 1.027026591533215
$inn[13,1] = 1.018060454885414
$inn[13,2] = 1.022450525842922
$inn[13,3] = 0.9921799884222134
$inn[13,4] = 1.026083433675736
$inn[13,5] = 1.010170496643154
$bf[14,0] = 1.02
$bf[14,1] = 1.011854398539234
$bf[14,2] = 1.019197643331845
$bf[14,3] = 0.9890781214508737
$bf[14,4] = 1.023148894182646
$inn[14,0] = 1.027119196494767
$inn[14,1] = 1.018341756156396
$inn[14,2] = 1.022696130074645
$inn[14,3] = 0.9926916645766087
$inn[14,4] = 1.026632843175832
$inn[14,5] = 1.010263132097421
$bf[15,0] = 1.02
$bf[15,1] = 1.012172015349702
$bf[15,2] = 1.019425242260828
$bf[15,3] = 0.989476357848556
$bf[15,4] = 1.023567578900809
$inn[15,0] = 1.027176687485335
$inn[15,1] = 1.0185180533195
$inn[15,2] = 1.022849882188251
$inn[15,3] = 0.9930127773699352
$inn[15,4] = 1.026977326816189
$inn[15,5] = 1.010321186129477
$bf[16,0] = 1.02
$bf[16,1] = 1.012357294893615
$bf[16,2] = 1.019557968800318
$bf[16,3] = 0.9897087662937556
$bf[16,4] = 1.023811806712793
$inn[16,0] = 1.027210005393617
$inn[16,1] = 1.018620828056858
$inn[16,2] = 1.022939451507028
$inn[16,3] = 0.9932001317071769
$inn[16,4] = 1.027178205330314
$inn[16,5] = 1.010355028591221
$bf[17,0] = 1.02
$bf[17,1] = 1.012420473620612
$bf[17,2] = 1.019603220311462
$bf[17,3] = 0.9897880325774034
$bf[17,4] = 1.023895084800613
$inn[17,0] = 1.027221329345908
$inn[17,1] = 1.018655861988369
$inn[17,2] = 1.022969973352381
$inn[17,3] = 0.9932640239640975
$inn[17,4] = 1.027246690710084
$inn[17,5] = 1.010366564683159
$bf[18,0] = 1.02
$bf[18,1] = 1.012137936094309
$bf[18,2] = 1.019400825941109
$bf[18,3] = 0.9894336180360679
$bf[18,4] = 1.023522656323931
$inn[18,0] = 1.027170541542954
$inn[18,1] = 1.018499144130744
$inn[18,2] = 1.022833397585991
$inn[18,3] = 0.9929783193494215
$inn[18,4] = 1.026940372455662
$inn[18,5] = 1.010314959493056
$bf[19,0] = 1.02
$bf[19,1] = 1.011219605747723
$bf[19,2] = 1.018742496003814
$bf[19,3] = 0.9882828385668249
$bf[19,4] = 1.02231204569826
$inn[19,0] = 1.027002892915203
$inn[19,1] = 1.01798897830507
$inn[19,2] = 1.022388066343547
$inn[19,3] = 0.9920501090198102
$inn[19,4] = 1.02594388105627
$inn[19,5] = 1.010146957899019
$bf[20,0] = 1.02
$bf[20,1] = 1.010642394957656
$bf[20,2] = 1.018328335941057
$bf[20,3] = 0.9875604150241495
$bf[20,4] = 1.021551035259046
$inn[20,0] = 1.026895547968592
$inn[20,1] = 1.017667715097555
$inn[20,2] = 1.022107070267467
$inn[20,3] = 0.9914670000341481
$inn[20,4] = 1.025316874388953
$inn[20,5] = 1.010041155451712
$bf[21,0] = 1.02
$bf[21,1] = 1.010948367850734
$bf[21,2] = 1.018547912427679
$bf[21,3] = 0.9879432794643023
$bf[21,4] = 1.021954447046985
$inn[21,0] = 1.02695263698338
$inn[21,1] = 1.017838070042108
$inn[21,2] = 1.022256125900024
$inn[21,3] = 0.991776070289318
$inn[21,4] = 1.025649307116903
$inn[21,5] = 1.010097259637122
$bf[22,0] = 1.02
$bf[22,1] = 1.012153334984196
$bf[22,2] = 1.01941185871189
$bf[22,3] = 0.9894529299347244
$bf[22,4] = 1.023542954852929
$inn[22,0] = 1.027173319296033
$inn[22,1] = 1.018507688553604
$inn[22,2] = 1.02284084661338
$inn[22,3] = 0.9929938892766442
$inn[22,4] = 1.026957070704106
$inn[22,5] = 1.01031777310223
$bf[23,0] = 1.02
$bf[23,1] = 1.013552374324425
$bf[23,2] = 1.020413313496076
$bf[23,3] = 0.9912096547607049
$bf[23,4] = 1.025386959339229
$inn[23,0] = 1.027420948107964
$inn[23,1] = 1.019282527000141
$inn[23,2] = 1.023514997110262
$inn[23,3] = 0.9944092447426414
$inn[23,4] = 1.028472577516152
$inn[23,5] = 1.010572901609552

$ws.Range("B2:F25").Value = $bf
$ws.Range("I2:N25").Value = $inn
